$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old content (old template had 5 columns: Province Code, English Name, Local Name, Country, Inactive/Active)
$ws.Range("A1:E2").Clear()

# Headers
$ws.Range("A1").Value = "Province Code"
$ws.Range("B1").Value = "English Name"
$ws.Range("C1").Value = "Local Name"
$ws.Range("D1").Value = "Country Code"
$ws.Range("E1").Value = "Postal Code"
$ws.Range("F1").Value = "Inactive"

# Values
$ws.Range("A2").Value = "HAN"
$ws.Range("B2").Value = " Ha Noi"
$ws.Range("C2").Value = "Hà Nội"
$ws.Range("D2").Value = "VN"
$ws.Range("E2").Value = 100000
$ws.Range("F2").Value = "Active"

$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1:F1").Borders.LineStyle = 1
$ws.Range("A1:F1").Borders.Weight = 2

$ws.Columns.Item(1).ColumnWidth = 18
$ws.Columns.Item(2).ColumnWidth = 15.85546875
$ws.Columns.Item(3).ColumnWidth = 13.85546875
$ws.Columns.Item(4).ColumnWidth = 20.85546875
$ws.Columns.Item(5).ColumnWidth = 14.42578125
